$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the underlying numbers; the Recovery rate (E) column formulas
# will recalculate automatically since E = D/C
$ws.Range("C2").Value = 1936832994
$ws.Range("D2").Value = 547732842

$ws.Range("C3").Value = 7848912935
$ws.Range("D3").Value = 1971679882

$ws.Range("C4").Value = 4593972824
$ws.Range("D4").Value = 1080293151

$ws.Range("C5").Value = 6887766656
$ws.Range("D5").Value = 1518357404

$ws.Range("C6").Value = 7807949106
$ws.Range("D6").Value = 1718944460

# Update the active selection to match the post-edit state
$ws.Range("E2:E6").Select()
